$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New warm fall schedule: flatten the midday Temperature (column C) plateau
# from 21/22 down to a steady 20 for rows 10-19.
$ws.Range("C10").Value = 20
$ws.Range("C11").Value = 20
$ws.Range("C12").Value = 20
$ws.Range("C13").Value = 20
$ws.Range("C14").Value = 20
$ws.Range("C15").Value = 20
$ws.Range("C16").Value = 20
$ws.Range("C17").Value = 20
$ws.Range("C18").Value = 20
$ws.Range("C19").Value = 20

# Update the active selection to match where the editor left off.
$ws.Range("H32").Select()
